$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 446.23254
$ws.Range("J17").Value = 454.7619
$ws.Range("L17").Value = 1364.2857
$ws.Range("N17").Value = -1700.2857
$ws.Range("H111").Value = 1193.5834
$ws.Range("I111").Value = 1079.8572
$ws.Range("J111").Value = 1352.8
$ws.Range("K111").Value = 3239.5716
$ws.Range("L111").Value = 4058.4
$ws.Range("M111").Value = -172.5715999999998
$ws.Range("N111").Value = -10192.4
$ws.Range("H113").Value = 2086.2727
$ws.Range("I113").Value = 1775
$ws.Range("J113").Value = 2916.3333
$ws.Range("K113").Value = 1775
$ws.Range("L113").Value = 2916.3333
$ws.Range("M113").Value = 1479
$ws.Range("N113").Value = -9424.3333
$ws.Range("H127").Value = 1638.7273
$ws.Range("I127").Value = 671.375
$ws.Range("J127").Value = 4218.3335
$ws.Range("K127").Value = 2014.125
$ws.Range("L127").Value = 12655.0005
$ws.Range("M127").Value = 2945.875
$ws.Range("N127").Value = -22575.0005
$ws.Range("H129").Value = 19451.777
$ws.Range("J129").Value = 20956.38
$ws.Range("L129").Value = 62869.14
$ws.Range("N129").Value = -72869.14
$ws.Range("H132").Value = 4465405
$ws.Range("I132").Value = 5103009
$ws.Range("J132").Value = 2176.5
$ws.Range("K132").Value = 15309027
$ws.Range("L132").Value = 6529.5
$ws.Range("M132").Value = -15306497
$ws.Range("N132").Value = -11589.5
$ws.Range("H138").Value = 1269.14
$ws.Range("I138").Value = 857.34283
$ws.Range("K138").Value = 2572.02849
$ws.Range("M138").Value = 2567.97151

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3336.068
$ws.Range("I32").Value = 3087.2703
$ws.Range("K32").Value = 3087.2703
$ws.Range("M32").Value = -2800.2703
$ws.Range("H45").Value = 931.7619
$ws.Range("I45").Value = 980.5
$ws.Range("J45").Value = 901.7692
$ws.Range("K45").Value = 980.5
$ws.Range("L45").Value = 901.7692
$ws.Range("M45").Value = -603.5
$ws.Range("N45").Value = -1655.7692
$ws.Range("H61").Value = 2530.2354
$ws.Range("I61").Value = 2042.8334
$ws.Range("K61").Value = 2042.8334
$ws.Range("M61").Value = -1830.8334
$ws.Range("H132").Value = 6621.913
$ws.Range("I132").Value = 7606.1763
$ws.Range("J132").Value = 3833.1667
$ws.Range("K132").Value = 22818.5289
$ws.Range("L132").Value = 11499.5001
$ws.Range("M132").Value = -20288.5289
$ws.Range("N132").Value = -16559.5001
$ws.Range("H136").Value = 2530.2354
$ws.Range("I136").Value = 2042.8334
$ws.Range("K136").Value = 6128.5002
$ws.Range("M136").Value = -3578.5002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1054.7
$ws.Range("I58").Value = 1123.5385
$ws.Range("J58").Value = 926.8570999999999
$ws.Range("K58").Value = 1123.5385
$ws.Range("L58").Value = 926.8570999999999
$ws.Range("M58").Value = -920.5385000000001
$ws.Range("N58").Value = -1332.8571
$ws.Range("H99").Value = 2352.5881
$ws.Range("I99").Value = 1666.6666
$ws.Range("K99").Value = 1666.6666
$ws.Range("M99").Value = -168.6666
$ws.Range("H126").Value = 2352.5881
$ws.Range("I126").Value = 1666.6666
$ws.Range("K126").Value = 4999.9998
$ws.Range("M126").Value = -2529.9998
$ws.Range("H136").Value = 1054.7
$ws.Range("I136").Value = 1123.5385
$ws.Range("J136").Value = 926.8570999999999
$ws.Range("K136").Value = 3370.6155
$ws.Range("L136").Value = 2780.5713
$ws.Range("M136").Value = -820.6155000000003
$ws.Range("N136").Value = -7880.5713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 409.57144
$ws.Range("I2").Value = 471.16666
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 2826.99996
$ws.Range("L2").Value = 240
$ws.Range("M2").Value = -2713.99996
$ws.Range("N2").Value = -466
$ws.Range("H12").Value = 73.64706
$ws.Range("I12").Value = 3.3333333
$ws.Range("J12").Value = 88.71429000000001
$ws.Range("K12").Value = 9.999999900000001
$ws.Range("L12").Value = 266.14287
$ws.Range("M12").Value = 163.0000001
$ws.Range("N12").Value = -612.14287
$ws.Range("H63").Value = 9615.143
$ws.Range("I63").Value = 2903.75
$ws.Range("J63").Value = 11603.704
$ws.Range("K63").Value = 8711.25
$ws.Range("L63").Value = 34811.112
$ws.Range("M63").Value = -7962.25
$ws.Range("N63").Value = -36309.112
$ws.Range("H66").Value = 9615.143
$ws.Range("I66").Value = 2903.75
$ws.Range("J66").Value = 11603.704
$ws.Range("K66").Value = 26133.75
$ws.Range("L66").Value = 104433.336
$ws.Range("M66").Value = -22389.75
$ws.Range("N66").Value = -111921.336
$ws.Range("H113").Value = 502.77585
$ws.Range("I113").Value = 500.6875
$ws.Range("J113").Value = 505.34616
$ws.Range("K113").Value = 1502.0625
$ws.Range("L113").Value = 1516.03848
$ws.Range("M113").Value = 667.9375
$ws.Range("N113").Value = -5856.03848
$ws.Range("H131").Value = 3620544.8
$ws.Range("J131").Value = 8624651
$ws.Range("L131").Value = 25873953
$ws.Range("N131").Value = -25884033
$ws.Range("H137").Value = 39696812
$ws.Range("I137").Value = 47619736
$ws.Range("J137").Value = 35735348
$ws.Range("K137").Value = 142859208
$ws.Range("L137").Value = 107206044
$ws.Range("M137").Value = -142854108
$ws.Range("N137").Value = -107216244
$ws.Range("H141").Value = 1915.4546
$ws.Range("I141").Value = 1748.5186
$ws.Range("J141").Value = 2666.6667
$ws.Range("K141").Value = 5245.5558
$ws.Range("L141").Value = 8000.000100000001
$ws.Range("M141").Value = -65.55580000000009
$ws.Range("N141").Value = -18360.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15943248
$ws.Range("I70").Value = 26846404
$ws.Range("J70").Value = 7868.4614
$ws.Range("K70").Value = 26846404
$ws.Range("L70").Value = 7868.4614
$ws.Range("M70").Value = -26846134
$ws.Range("N70").Value = -8408.4614
$ws.Range("H73").Value = 15943248
$ws.Range("I73").Value = 26846404
$ws.Range("J73").Value = 7868.4614
$ws.Range("K73").Value = 26846404
$ws.Range("L73").Value = 7868.4614
$ws.Range("M73").Value = -26845468
$ws.Range("N73").Value = -9740.4614
$ws.Range("H102").Value = 1792.091
$ws.Range("I102").Value = 1756.8572
$ws.Range("J102").Value = 1853.75
$ws.Range("K102").Value = 1756.8572
$ws.Range("L102").Value = 1853.75
$ws.Range("M102").Value = -134.8571999999999
$ws.Range("N102").Value = -5097.75
$ws.Range("H107").Value = 717.6316
$ws.Range("I107").Value = 642.3333
$ws.Range("K107").Value = 642.3333
$ws.Range("M107").Value = 1277.6667
$ws.Range("H113").Value = 20834466
$ws.Range("I113").Value = 41667430
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 41667430
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -41665260
$ws.Range("N113").Value = -5840
$ws.Range("H132").Value = 103640.1
$ws.Range("I132").Value = 121029.586
$ws.Range("K132").Value = 363088.758
$ws.Range("M132").Value = -360558.758

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2497.6365
$ws.Range("I40").Value = 2497.6365
$ws.Range("K40").Value = 2497.6365
$ws.Range("M40").Value = -2361.6365
$ws.Range("H61").Value = 1761
$ws.Range("I61").Value = 1100
$ws.Range("J61").Value = 2752.5
$ws.Range("K61").Value = 1100
$ws.Range("L61").Value = 2752.5
$ws.Range("M61").Value = -898
$ws.Range("N61").Value = -3156.5
$ws.Range("H113").Value = 1761
$ws.Range("I113").Value = 1100
$ws.Range("J113").Value = 2752.5
$ws.Range("K113").Value = 1100
$ws.Range("L113").Value = 2752.5
$ws.Range("M113").Value = 1070
$ws.Range("N113").Value = -7092.5
$ws.Range("H122").Value = 4250.5
$ws.Range("I122").Value = 5501
$ws.Range("K122").Value = 16503
$ws.Range("M122").Value = -14053

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 12500
$ws.Range("I26").Value = 15000
$ws.Range("J26").Value = 10000
$ws.Range("K26").Value = 15000
$ws.Range("L26").Value = 10000
$ws.Range("M26").Value = -14707
$ws.Range("N26").Value = -10586
$ws.Range("H107").Value = 621.4
$ws.Range("I107").Value = 651.3333
$ws.Range("J107").Value = 576.5
$ws.Range("K107").Value = 1953.9999
$ws.Range("L107").Value = 1729.5
$ws.Range("M107").Value = -33.99990000000003
$ws.Range("N107").Value = -5569.5
$ws.Range("H113").Value = 522.63635
$ws.Range("I113").Value = 484.9
$ws.Range("K113").Value = 1454.7
$ws.Range("M113").Value = 715.3000000000002
$ws.Range("H122").Value = 1501.4286
$ws.Range("I122").Value = 1580.7693
$ws.Range("J122").Value = 1372.5
$ws.Range("K122").Value = 4742.3079
$ws.Range("L122").Value = 1372.5
$ws.Range("M122").Value = -2292.3079
$ws.Range("N122").Value = -9017.5
